# Auto-generated edit script: adds 21 new trivia rows to Hoja1 (sheet1),
# 2 new rows to Hoja2 (sheet2), and scrolls Hoja3 (sheet3)'s view.
# Commit message: 'Changed activity transitions into within-activity transitions'

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")
$ws3 = $wb.Worksheets.Item("Hoja3")

# ---------------------------------------------------------------------
# Cell writes, issued in the same order the source data was authored in
# (this matters only for the internal shared-string table ordering; the
# resulting cell values are identical regardless of order).
# ---------------------------------------------------------------------
$writes = @(
    @{ Ws = $ws2; Row = 69; Col = 'B'; Value = 'John Corzine' }
    @{ Ws = $ws2; Row = 70; Col = 'B'; Value = 'Hunt Brothers' }
    @{ Ws = $ws1; Row = 137; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 138; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 139; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 140; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 141; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 142; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 143; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 144; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 145; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 146; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 147; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 148; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 149; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 150; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 151; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 152; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 153; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 154; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 155; Col = 'A'; Value = 10 }
    @{ Ws = $ws1; Row = 137; Col = 'B'; Value = 'Bernie Madoff performed duties as. . .' }
    @{ Ws = $ws1; Row = 137; Col = 'D'; Value = 'NASDAQ CEO' }
    @{ Ws = $ws1; Row = 137; Col = 'E'; Value = 'White-collar stealing' }
    @{ Ws = $ws1; Row = 137; Col = 'F'; Value = 'Stock analyst' }
    @{ Ws = $ws1; Row = 138; Col = 'D'; Value = 'Bernie Madoff' }
    @{ Ws = $ws1; Row = 138; Col = 'B'; Value = 'Who was the head of a pyramidal scheme that notoriously blew up in 2011?' }
    @{ Ws = $ws1; Row = 138; Col = 'E'; Value = 'Lloyd Blankfein' }
    @{ Ws = $ws1; Row = 138; Col = 'F'; Value = 'Blythe Masters' }
    @{ Ws = $ws1; Row = 139; Col = 'D'; Value = 'WTI futures' }
    @{ Ws = $ws1; Row = 139; Col = 'E'; Value = 'Brent futures' }
    @{ Ws = $ws1; Row = 139; Col = 'F'; Value = 'Cushing spot swaps' }
    @{ Ws = $ws1; Row = 139; Col = 'B'; Value = 'What contracts are traded at NYMEX?' }
    @{ Ws = $ws1; Row = 140; Col = 'B'; Value = 'Where does the Brent future contract trade?' }
    @{ Ws = $ws1; Row = 140; Col = 'D'; Value = 'ICE' }
    @{ Ws = $ws1; Row = 140; Col = 'E'; Value = 'NYMEX' }
    @{ Ws = $ws1; Row = 140; Col = 'F'; Value = 'CBOE' }
    @{ Ws = $ws1; Row = 141; Col = 'B'; Value = 'What does Brent oil mean?' }
    @{ Ws = $ws1; Row = 141; Col = 'D'; Value = 'Oil pumped from the Norwegian North Sea coast' }
    @{ Ws = $ws1; Row = 141; Col = 'E'; Value = 'Oil mixed with heavy sands' }
    @{ Ws = $ws1; Row = 141; Col = 'F'; Value = 'Oil from the OPEC' }
    @{ Ws = $ws1; Row = 142; Col = 'B'; Value = 'What is the oil that sets gas prices in Europe?' }
    @{ Ws = $ws1; Row = 142; Col = 'D'; Value = 'Brent oil' }
    @{ Ws = $ws1; Row = 142; Col = 'E'; Value = 'OPEC oil' }
    @{ Ws = $ws1; Row = 142; Col = 'F'; Value = 'WTI oil' }
    @{ Ws = $ws1; Row = 143; Col = 'B'; Value = 'Brent oil contracts are traded. . .' }
    @{ Ws = $ws1; Row = 143; Col = 'D'; Value = 'In London and priced in dollars' }
    @{ Ws = $ws1; Row = 143; Col = 'E'; Value = 'In London and priced in pounds sterling' }
    @{ Ws = $ws1; Row = 143; Col = 'F'; Value = 'In New York City and priced in dollars' }
    @{ Ws = $ws1; Row = 144; Col = 'B'; Value = 'What can I pick up in Cushing, Ocklahoma?' }
    @{ Ws = $ws1; Row = 144; Col = 'D'; Value = 'WTI barrels' }
    @{ Ws = $ws1; Row = 144; Col = 'E'; Value = 'Live cattle' }
    @{ Ws = $ws1; Row = 144; Col = 'F'; Value = 'Wheat bushels' }
    @{ Ws = $ws1; Row = 145; Col = 'B'; Value = 'Where is the strait of Hormuz?' }
    @{ Ws = $ws1; Row = 145; Col = 'D'; Value = 'In the Persial gulf' }
    @{ Ws = $ws1; Row = 145; Col = 'E'; Value = 'In Southern Spain' }
    @{ Ws = $ws1; Row = 146; Col = 'B'; Value = 'What is the OPEC?' }
    @{ Ws = $ws1; Row = 146; Col = 'D'; Value = 'An oil cartel composed mostly of Arab nations with headquartes in Vienna' }
    @{ Ws = $ws1; Row = 146; Col = 'E'; Value = 'An oil company' }
    @{ Ws = $ws1; Row = 146; Col = 'F'; Value = 'An oil cartel with headquartes in Qatar' }
    @{ Ws = $ws1; Row = 145; Col = 'F'; Value = 'In Northern Iran' }
    @{ Ws = $ws1; Row = 149; Col = 'E'; Value = 'I believe that banking institutions are more dangerous to our liberties than standing armies' }
    @{ Ws = $ws1; Row = 147; Col = 'B'; Value = 'Who said "I believe that banking institutions are more dangerous to our liberties than standing armies"?' }
    @{ Ws = $ws1; Row = 147; Col = 'D'; Value = 'Thomas Jefferson' }
    @{ Ws = $ws1; Row = 147; Col = 'E'; Value = 'Benjamin Franklin' }
    @{ Ws = $ws1; Row = 147; Col = 'F'; Value = 'John Adams' }
    @{ Ws = $ws1; Row = 148; Col = 'B'; Value = 'Who said "The modern theory of the perpetuation of debt has drenched the earth with blood, and crushed its inhabitants under burdens ever accumulating"?' }
    @{ Ws = $ws1; Row = 148; Col = 'D'; Value = 'Thomas Jefferson' }
    @{ Ws = $ws1; Row = 148; Col = 'E'; Value = 'James Madison' }
    @{ Ws = $ws1; Row = 148; Col = 'F'; Value = 'Andrew Jackson' }
    @{ Ws = $ws1; Row = 149; Col = 'B'; Value = 'In his speech against the 1809 recharter of the First Bank, Thomas Jefferson said. . .' }
    @{ Ws = $ws1; Row = 149; Col = 'D'; Value = '. . .The issuing power should be taken from the banks and restored to the people, to whom it properly belongs' }
    @{ Ws = $ws1; Row = 149; Col = 'F'; Value = '"Congress was given the reight to issue paper money themselves, not to delegate it to individuals or corporations"' }
    @{ Ws = $ws1; Row = 150; Col = 'B'; Value = 'Who said "Issue of currency should be lodged with the government and be protected from domination by Wall Street"?' }
    @{ Ws = $ws1; Row = 150; Col = 'D'; Value = 'Theodore Roosevelt' }
    @{ Ws = $ws1; Row = 150; Col = 'E'; Value = 'James Madison' }
    @{ Ws = $ws1; Row = 150; Col = 'F'; Value = 'Benjamin Franklin' }
    @{ Ws = $ws1; Row = 151; Col = 'B'; Value = 'Who said "If congress has the right under the Constitution to issue paper money, it was  given them to use themselves, not to be delegated to individuals or corporations"?' }
    @{ Ws = $ws1; Row = 151; Col = 'D'; Value = 'Andrew Jackson' }
    @{ Ws = $ws1; Row = 151; Col = 'E'; Value = 'James Madison' }
    @{ Ws = $ws1; Row = 151; Col = 'F'; Value = 'Thomas Jefferson' }
    @{ Ws = $ws1; Row = 152; Col = 'B'; Value = 'The words "I am a most unhappy man. I have unwittingly ruined my country…" were written by. . .' }
    @{ Ws = $ws1; Row = 152; Col = 'D'; Value = 'Woodrow Wilson' }
    @{ Ws = $ws1; Row = 152; Col = 'E'; Value = 'Theodor Roosevelt' }
    @{ Ws = $ws1; Row = 152; Col = 'F'; Value = 'Franklin D. Roosevelt' }
    @{ Ws = $ws1; Row = 153; Col = 'B'; Value = 'Who said "The real truth of the matter is that a financial element in the large centers has owned the government ever since the days of Andrew Jackson"?' }
    @{ Ws = $ws1; Row = 153; Col = 'D'; Value = 'Franklin D. Roosevelt' }
    @{ Ws = $ws1; Row = 153; Col = 'E'; Value = 'Woodrow Wilson' }
    @{ Ws = $ws1; Row = 153; Col = 'F'; Value = 'Theodor Roosevelt' }
    @{ Ws = $ws1; Row = 154; Col = 'B'; Value = 'Who said "Money has no motherland; financiers are without patriotism and without decency; their sole object is gain"?' }
    @{ Ws = $ws1; Row = 154; Col = 'D'; Value = 'Napoleon Bonaparte' }
    @{ Ws = $ws1; Row = 154; Col = 'E'; Value = 'Thomas Jefferson' }
    @{ Ws = $ws1; Row = 154; Col = 'F'; Value = 'Benjamin Franklin' }
    @{ Ws = $ws1; Row = 155; Col = 'B'; Value = 'Who said "I fear that foreign bankers will entirely control the exuberant riches of America and use them to systematically corrupt civilization"?' }
    @{ Ws = $ws1; Row = 155; Col = 'D'; Value = 'Otto von Bismarck' }
    @{ Ws = $ws1; Row = 155; Col = 'E'; Value = 'Napoleon Bonaparte' }
    @{ Ws = $ws1; Row = 155; Col = 'F'; Value = 'Andrew Jackson' }
    @{ Ws = $ws1; Row = 156; Col = 'B'; Value = 'Where can you read the following text? "Money plays the largest part in determining the course of history"' }
    @{ Ws = $ws1; Row = 156; Col = 'D'; Value = 'The Communist Manifesto' }
    @{ Ws = $ws1; Row = 156; Col = 'E'; Value = 'The Republic' }
    @{ Ws = $ws1; Row = 156; Col = 'F'; Value = 'The Wealth of Nations' }
    @{ Ws = $ws1; Row = 157; Col = 'B'; Value = 'Who said "Banks lend by creating credit. They create the means of payment out of nothing"?' }
    @{ Ws = $ws1; Row = 157; Col = 'D'; Value = 'Ralph M. Hawtry' }
    @{ Ws = $ws1; Row = 157; Col = 'E'; Value = 'Otto von Bismarck' }
    @{ Ws = $ws1; Row = 157; Col = 'F'; Value = 'Napoleon Bonaparte' }
)

$colIndex = @{ A = 1; B = 2; C = 3; D = 4; E = 5; F = 6 }
foreach ($w in $writes) {
    $c = $colIndex[$w.Col]
    $w.Ws.Cells.Item($w.Row, $c).Value = $w.Value
}

# ---------------------------------------------------------------------
# Sheet view / selection updates
# ---------------------------------------------------------------------
$ws1.Application.ActiveWindow.ScrollRow = 142
$ws1.Range("A158").Select()

$ws2.Application.ActiveWindow.ScrollRow = 46
$ws2.Range("B71").Select()

$ws3.Application.ActiveWindow.ScrollRow = 7
$ws3.Range("D15").Select()

